# Rename Sheet1 -> data, add a new "coords" sheet with x/y data,
# and update selections/active tab as per the commit.

$wb = $excel.ActiveWorkbook

# Rename existing sheet to "data"
$dataSheet = $wb.Worksheets.Item(1)
$dataSheet.Name = "data"

# Update row 5 height on data sheet back to default (12.8)
$dataSheet.Rows.Item(5).RowHeight = 12.8

# Update selection on data sheet (no longer the active/selected tab)
$dataSheet.Range("C25").Select()

# Add a new worksheet named "coords" after "data"
$coordsSheet = $wb.Worksheets.Add($null, $dataSheet)
$coordsSheet.Name = "coords"

# Fill in header row
$coordsSheet.Range("B1").Value = "x"
$coordsSheet.Range("C1").Value = "y"

# Fill in well names and coordinates
$wells = @("Well_1", "Well_2", "Well_3", "Well_4", "Well_5", "Well_6", "Well_7")
$xs = @(30, 90, 70, 50, 60, 31.72, 60)
$ys = @(40, 70, 70, 90, 10, 73.28, 95)

for ($i = 0; $i -lt $wells.Length; $i++) {
    $row = $i + 2
    $coordsSheet.Range("A$row").Value = $wells[$i]
    $coordsSheet.Range("B$row").Value = $xs[$i]
    $coordsSheet.Range("C$row").Value = $ys[$i]
}

# Numeric coordinate cells are center-aligned in the source sheet
$coordsSheet.Range("B2:C8").HorizontalAlignment = -4108

# Select D14 on coords sheet (it becomes the active tab)
$coordsSheet.Range("D14").Select()

# coords sheet is now the active tab
$coordsSheet.Activate()

$wb.Save()
